$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from H1 (existing header style) into I1:J1 so the new
# headers match the look of the other header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I2:J63 values (I0, IF) for each data row, in row order starting at row 2.
$data = @(@(9,9),@(9,9),@(9,9),@(8,8),@(9,9),@(7,8),@(9,10),@(7,8),@(7,7),@(7,7),@(10,10),@(5,5),@(8,8),@(7,7),@(8,8),@(8,8),@(8,8),@(8,8),@(8,8),@(8,8),@(9,9),@(8,8),@(8,8),@(9,9),@(8,9),@(9,9),@(8,8),@(8,8),@(9,9),@(10,10),@(8,9),@(9,9),@(8,9),@(9,10),@(8,8),@(8,8),@(8,8),@(7,7),@(8,8),@(6,7),@(8,8),@(9,9),@(9,9),@(7,8),@(9,9),@(9,9),@(8,8),@(8,9),@(7,8),@(9,9),@(6,7),@(8,9),@(7,7),@(7,7),@(7,8),@(9,9),@(5,6),@(5,5),@(7,7),@(8,8),@(8,9),@(8,9))

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $r = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
